# Slide 13, "Content Placeholder 2" shape: the GitHub repo link run currently
# reads "https://github.com/rod53/PresentationSDID/main/Presentation_SDID.ipynb"
# followed by a blank trailing paragraph. Trim the displayed link text down to
# "https://github.com/rod53/PresentationSDID", drop the legacy hyperlink-color
# extension that rides along with the hlinkClick, and merge away the now
# pointless empty second paragraph.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The hyperlinked run is the first 70 characters of the text frame (the
# trailing character is the paragraph mark that starts the empty 2nd para).
$run = $tr.Characters(1, 70)

# Re-assert the hyperlink address on the run; this rewrites a clean
# <a:hlinkClick r:id="..."/> and drops the old ahyp:hlinkClr extLst blob.
$actionSetting = $run.ActionSettings.Item(1)
$hyperlink = $actionSetting.Hyperlink
$hyperlink.Address = "https://github.com/rod53/PresentationSDID/main/Presentation_SDID.ipynb"

# Shorten the visible link text.
$run.Text = "https://github.com/rod53/PresentationSDID"

# Remove the now-empty trailing paragraph, merging its endParaRPr onto the
# (now shorter) first paragraph.
$trAfter = $sh.TextFrame.TextRange
$secondParagraph = $trAfter.Paragraphs(2, 1)
$secondParagraph.Delete()
